$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header in row 1, data starts row 2).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($e -eq $null -or $f -eq $null) {
        continue
    }

    # F holds the start date encoded as an 8-digit YYYYMMDD number.
    # Rows whose F value isn't a well-formed date (e.g. malformed/typo'd
    # entries) are left untouched, matching the observed source behaviour.
    $fStr = [string]([long]$f)
    if ($fStr.Length -ne 8) {
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startDate = Get-Date -Year $year -Month $month -Day $day

    if ($e -eq 1) {
        # Remaining days hit zero -> cycle restarts: reset remaining to 10
        # and push the start date forward by 10 days.
        $newE = 10
        $newDate = $startDate.AddDays(10)
        $newF = [int]$newDate.ToString("yyyyMMdd")
    } else {
        # Otherwise simply count down one more day.
        $newE = $e - 1
        $newF = [int]$f
    }

    $eCell.Value = $newE
    $fCell.Value = $newF
}
